# Generate Report for Handoff
#
# - assigns a fresh GUID-based file name to the two e2e test rows
# - flips their status from "Handed back: in sync with en-US" to "Ready for handoff"
# - refreshes the Latest HO Xliff Generate Date / per-language handoff timestamps
# - marks the second row as a content duplicate of the first (zh-cn/de-de sheets)
# - resets the "Latest Handback" columns back to their not-yet-handed-back sentinel
#   and drops the now-meaningless "Latest Target File" hyperlink/value

$wb = $excel.ActiveWorkbook

$NEW1 = "d29aeb6d-c59e-4536-9793-bcd535208054"
$NEW2 = "ffff1c900cbd-ee74-48de-b092-3167c93c63ac"

$zhXlf = "$NEW1.5afe5d86e2ef23868f0970ede3cfccfcf7ee92ef.zh-cn.xlf"
$deXlf = "$NEW1.5afe5d86e2ef23868f0970ede3cfccfcf7ee92ef.de-de.xlf"

$statusText = "Ready for handoff"
$hoDate     = "2016-09-07 05:21:53"
$zhDate     = "2016-09-07 05:21:46"
$deDate     = "2016-09-07 05:21:53"
$nullDate   = "0001-01-01 00:00:00"

function Set-HyperlinkDisplay($worksheet, $addrTarget, $newDisplay) {
    foreach ($h in $worksheet.Hyperlinks) {
        if ($h.Range.Address() -eq $addrTarget) {
            $h.TextToDisplay = $newDisplay
            return
        }
    }
}

function Remove-HyperlinkAt($worksheet, $addrTarget) {
    foreach ($h in $worksheet.Hyperlinks) {
        if ($h.Range.Address() -eq $addrTarget) {
            $h.Delete()
            return
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-HyperlinkDisplay $wsOverview '$B$2' "e2e\$NEW1.md"
Set-HyperlinkDisplay $wsOverview '$B$3' "e2e\$NEW2.md"

$wsOverview.Range("A2").Value = "$NEW1.md"
$wsOverview.Range("B2").Value = "e2e\$NEW1.md"
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("G2").Value = $hoDate

$wsOverview.Range("A3").Value = "$NEW2.md"
$wsOverview.Range("B3").Value = "e2e\$NEW2.md"
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Range("G3").Value = $hoDate

$wsOverview.Range("E1").EntireColumn.ColumnWidth = 16.33
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 16.33

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-HyperlinkDisplay $wsZh '$A$2' "$NEW1.md"
Set-HyperlinkDisplay $wsZh '$A$3' "$NEW2.md"
Remove-HyperlinkAt $wsZh '$I$2'
Remove-HyperlinkAt $wsZh '$I$3'

$wsZh.Range("A2").Value = "$NEW1.md"
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("G2").Value = $zhXlf
$wsZh.Range("H2").Value = $zhDate
$wsZh.Range("I2").ClearFormats()
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").ClearFormats()
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $nullDate

$wsZh.Range("A3").Value = "$NEW2.md"
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $zhDate
$wsZh.Range("I3").ClearFormats()
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").ClearFormats()
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = $nullDate

$wsZh.Range("C1").EntireColumn.ColumnWidth = 16.33
$wsZh.Range("I1").EntireColumn.ColumnWidth = 17.83
$wsZh.Range("J1").EntireColumn.ColumnWidth = 20.83

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-HyperlinkDisplay $wsDe '$A$2' "$NEW1.md"
Set-HyperlinkDisplay $wsDe '$A$3' "$NEW2.md"
Remove-HyperlinkAt $wsDe '$I$2'
Remove-HyperlinkAt $wsDe '$I$3'

$wsDe.Range("A2").Value = "$NEW1.md"
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("G2").Value = $deXlf
$wsDe.Range("H2").Value = $deDate
$wsDe.Range("I2").ClearFormats()
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").ClearFormats()
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $nullDate

$wsDe.Range("A3").Value = "$NEW2.md"
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $deDate
$wsDe.Range("I3").ClearFormats()
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").ClearFormats()
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = $nullDate

$wsDe.Range("C1").EntireColumn.ColumnWidth = 16.33
$wsDe.Range("I1").EntireColumn.ColumnWidth = 17.83
$wsDe.Range("J1").EntireColumn.ColumnWidth = 20.83

Write-Output "Generate Report for Handoff: edits applied"
